# Applies the NATMI Fgf16-Fgfr4 update described in the commit "Natmi following Dr Hou advice".
# Rebuilds data rows 2-9 (sending/target cluster pairs ECs<->sCs with intermediate FAPs/M2 targets)
# using updated ligand-receptor specificity figures, and replaces the old single-row placeholder data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{}
$rowData[2] = @("ECs", "Fgf16", "Fgfr4", "ECs", 1, 0.3333333333333333, 0.096573, 0.289719, 0.1488791629603479, 0.2078467214814188, 2, 0.6666666666666666, 0.3341516666666666, 1.002455, 0.02354483703663662, 0.03467746452042096, 0.032270028905, 0.290430260145, 0.003505335630052258, 0.007207597309857716)
$rowData[3] = @("ECs", "Fgf16", "Fgfr4", "FAPs", 1, 0.3333333333333333, 0.096573, 0.289719, 0.1488791629603479, 0.2078467214814188, 1, 0.3333333333333333, 0.153566, 0.460698, 0.01082049501783563, 0.01593671391696275, 0.014830329318, 0.133472963862, 0.001610946241071984, 0.003312393738828008)
$rowData[4] = @("ECs", "Fgf16", "Fgfr4", "M2", 1, 0.3333333333333333, 0.096573, 0.289719, 0.1488791629603479, 0.2078467214814188, 1, 0.3333333333333333, 0.03596566666666667, 0.107897, 0.002534195830976933, 0.003732431270590561, 0.003473312327, 0.031259810943, 0.0003772889540934491, 0.0007757736027469742)
$rowData[5] = @("ECs", "Fgf16", "Fgfr4", "sCs", 1, 0.3333333333333333, 0.096573, 0.289719, 0.1488791629603479, 0.2078467214814188, 2, 1, 13.6684585, 27.336917, 0.9631004721145509, 0.9456533902920258, 1.3200040427205, 7.920024256323, 0.1433855921351302, 0.1965509568299861)
$rowData[6] = @("sCs", "Fgf16", "Fgfr4", "ECs", 2, 1, 0.5520940000000001, 1.104188, 0.8511208370396521, 0.7921532785185812, 2, 0.6666666666666666, 0.3341516666666666, 1.002455, 0.02354483703663662, 0.03467746452042096, 0.1844831302566667, 1.10689878154, 0.02003950140658436, 0.02746986721056324)
$rowData[7] = @("sCs", "Fgf16", "Fgfr4", "FAPs", 2, 1, 0.5520940000000001, 1.104188, 0.8511208370396521, 0.7921532785185812, 1, 0.3333333333333333, 0.153566, 0.460698, 0.01082049501783563, 0.01593671391696275, 0.08478286720400002, 0.508697203224, 0.00920954877676365, 0.01262432017813475)
$rowData[8] = @("sCs", "Fgf16", "Fgfr4", "M2", 2, 1, 0.5520940000000001, 1.104188, 0.8511208370396521, 0.7921532785185812, 1, 0.3333333333333333, 0.03596566666666667, 0.107897, 0.002534195830976933, 0.003732431270590561, 0.01985642877266667, 0.119138572636, 0.002156906876883484, 0.002956657667843587)
$rowData[9] = @("sCs", "Fgf16", "Fgfr4", "sCs", 2, 1, 0.5520940000000001, 1.104188, 0.8511208370396521, 0.7921532785185812, 2, 1, 13.6684585, 27.336917, 0.9631004721145509, 0.9456533902920258, 7.546273927099001, 30.18509570839601, 0.8197148799794207, 0.7491024334620396)

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}
